# Weekly refresh of the "Pera" (Packham's Triumph) price series.
#
# A new weekly observation is written into row 114 (Fecha/Precio
# mínimo/máximo/promedio/Precio $/Kg change; Calidad, Volumen, Unidad,
# Origen and Kg/unidad stay as they already were). To make room without
# losing any history, every following observation (rows 115-166) shifts
# down by exactly one row, and the row that falls off the bottom (the
# old row 166) is re-appended as the new last row, 167. Row 113 (above
# the new observation) and the purely descriptive columns A, B, C, E-K
# (Mercado, Producto, Variedad, etc. - identical for the whole block)
# are untouched throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shiftFirstRow = 115   # first row that gets pulled down from the row above it
$lastRow = 167         # brand-new row created at the bottom of the block

# Observation columns that shift down by one row for rows 115..167.
$shiftCols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D,L,M,N,O,P,Q,R,S,T

# Descriptive columns constant for the whole block; only needed to seed
# the brand-new row 167.
$constCols = @(1, 2, 3, 5, 6, 7, 8, 9, 10, 11)  # A,B,C,E,F,G,H,I,J,K

# Seed row 167's constant metadata from the row above (166) - identical
# throughout this product's block.
foreach ($c in $constCols) {
    $ws.Cells.Item($lastRow, $c).Value2 = $ws.Cells.Item($lastRow - 1, $c).Value2
}

# The "Fecha" column (D) carries a date number-format style; make sure the
# freshly-created last row picks it up too (new cells default to General).
$ws.Cells.Item($lastRow, 4).NumberFormat = $ws.Cells.Item($lastRow - 1, 4).NumberFormat

# Walk bottom-up (167 down to 115) so each source row is read before it
# gets overwritten.
for ($r = $lastRow; $r -ge $shiftFirstRow; $r--) {
    $src = $r - 1
    foreach ($c in $shiftCols) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($src, $c).Value2
    }
}

# Finally, write this week's new observation into row 114 (Calidad,
# Volumen, Unidad de comercialización, Origen and Kg/unidad keep their
# prior values; only the date and the three prices + $/Kg change).
$ws.Cells.Item(114, 4).Value2 = 44523    # Fecha
$ws.Cells.Item(114, 14).Value2 = 11000   # Precio mínimo
$ws.Cells.Item(114, 15).Value2 = 12000   # Precio máximo
$ws.Cells.Item(114, 16).Value2 = 11500   # Precio promedio ponderado
$ws.Cells.Item(114, 19).Value2 = 719     # Precio $/Kg
